# Helper: write a value into a cell while forcing it to be stored as TEXT
# (not auto-converted to a number), then restore the cell's original
# "General" number format / default style so no stray formatting is left
# behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "25.957.87"
Set-TextValue $ws.Range("E2") "  -0.75%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.745.91"
Set-TextValue $ws.Range("E3") "  -0.26%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  +0.08%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "248.76"
Set-TextValue $ws.Range("E5") "  +4.78%  "

# Row 6 - USDC
Set-TextValue $ws.Range("D6") "1.000"
Set-TextValue $ws.Range("E6") "  +0.11%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.5050"
Set-TextValue $ws.Range("E7") "  -8.99%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.2754"
Set-TextValue $ws.Range("E8") "  -2.89%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("E9") "  +0.00%  "

# Row 10 - now WrappedEther (was TRON)
Set-TextValue $ws.Range("B10") "WrappedEther"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D10") "1.747.94"
Set-TextValue $ws.Range("E10") "  -0.22%  "

# Row 11 - now TRON (was WrappedEther)
Set-TextValue $ws.Range("B11") "TRON"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D11") "0.07272"
Set-TextValue $ws.Range("E11") "  +0.82%  "

# Row 12 - Polygon
Set-TextValue $ws.Range("E12") "  +0.28%  "

# Row 13 - Solana
Set-TextValue $ws.Range("E13") "  -2.73%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("E14") "  +0.09%  "

# Row 15 - Litecoin
Set-TextValue $ws.Range("D15") "77.74"
Set-TextValue $ws.Range("E15") "  -1.10%  "

# Row 16 - Dai
Set-TextValue $ws.Range("E16") "  +0.11%  "

# Row 17 - BinanceUSD
Set-TextValue $ws.Range("E17") "  +0.09%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "25.975.41"
Set-TextValue $ws.Range("E18") "  -0.25%  "

# Row 19 - Avalanche
Set-TextValue $ws.Range("E19") "  +0.59%  "

# Row 20 - ShibaInu
Set-TextValue $ws.Range("D20") "0.000006846"
Set-TextValue $ws.Range("E20") "  +0.58%  "

# Row 21 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D21") "1.967.61"
Set-TextValue $ws.Range("E21") "  -0.14%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "4.457"
Set-TextValue $ws.Range("E22") "  +2.44%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("D23") "8.731"
Set-TextValue $ws.Range("E23") "  -0.27%  "

# Row 24 - Chainlink
Set-TextValue $ws.Range("D24") "5.402"
Set-TextValue $ws.Range("E24") "  +2.50%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "136.99"
Set-TextValue $ws.Range("E25") "  -1.92%  "

# Row 26 - Toncoin
Set-TextValue $ws.Range("D26") "1.508"
Set-TextValue $ws.Range("E26") "  -0.86%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("E27") "  -0.65%  "

# Row 28 - LidoDAOToken
Set-TextValue $ws.Range("D28") "1.782"
Set-TextValue $ws.Range("E28") "  -1.92%  "

# Row 29 - BitcoinCash
Set-TextValue $ws.Range("D29") "105.84"
Set-TextValue $ws.Range("E29") "  +0.16%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D30") "3.877"
Set-TextValue $ws.Range("E30") "  +1.33%  "

# Row 31 - Stellar
Set-TextValue $ws.Range("D31") "0.08201"
Set-TextValue $ws.Range("E31") "  -2.81%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "3.651"
Set-TextValue $ws.Range("E32") "  -0.20%  "

# Row 33 - Hedera
Set-TextValue $ws.Range("D33") "0.04686"
Set-TextValue $ws.Range("E33") "  +0.71%  "

# Row 34 - HuobiToken
Set-TextValue $ws.Range("D34") "2.653"
Set-TextValue $ws.Range("E34") "  +0.50%  "

# Row 35 - ARBITRUM
Set-TextValue $ws.Range("E35") "  -1.60%  "

# Row 36 - ImmutableX
Set-TextValue $ws.Range("D36") "0.6184"
Set-TextValue $ws.Range("E36") "  -2.42%  "

# Row 37 - MXToken
Set-TextValue $ws.Range("D37") "2.748"
Set-TextValue $ws.Range("E37") "  +1.35%  "

# Row 38 - VeChain
Set-TextValue $ws.Range("D38") "0.01615"
Set-TextValue $ws.Range("E38") "  -0.70%  "

# Row 39 - RenderToken
Set-TextValue $ws.Range("D39") "1.926"
Set-TextValue $ws.Range("E39") "  -3.17%  "

# Row 40 - PaxDollar
Set-TextValue $ws.Range("D40") "0.9998"
Set-TextValue $ws.Range("E40") "  +0.14%  "

# Row 41 - Quant
Set-TextValue $ws.Range("D41") "100.93"
Set-TextValue $ws.Range("E41") "  -1.40%  "

# Row 42 - TheSandbox
Set-TextValue $ws.Range("D42") "0.3934"
Set-TextValue $ws.Range("E42") "  -0.73%  "

# Row 43 - TrustWalletToken
Set-TextValue $ws.Range("D43") "0.7614"
Set-TextValue $ws.Range("E43") "  +1.69%  "

# Row 44 - FraxShare
Set-TextValue $ws.Range("E44") "  -1.89%  "

# Row 45 - Algorand
Set-TextValue $ws.Range("E45") "  -0.46%  "

# Row 46 - Aptos
Set-TextValue $ws.Range("D46") "6.337"
Set-TextValue $ws.Range("E46") "  -0.96%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "55.87"
Set-TextValue $ws.Range("E47") "  +1.82%  "

# Row 48 - Cronos
Set-TextValue $ws.Range("D48") "0.05299"
Set-TextValue $ws.Range("E48") "  -0.81%  "

# Row 49 - Elrond
Set-TextValue $ws.Range("D49") "30.73"
Set-TextValue $ws.Range("E49") "  -1.14%  "

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "7.595"
Set-TextValue $ws.Range("E50") "  -0.39%  "

# Row 51 - Decentraland
Set-TextValue $ws.Range("D51") "0.3444"
Set-TextValue $ws.Range("E51") "  -1.91%  "
